# Clue_LayoutStudent.xlsx edit script
# Commit message: "added new room and tests"
#
# This carves a new room ("J" / "JL" - e.g. a Jail/new room) out of the
# existing "M" room on Sheet1, re-colors a set of "room label" marker
# cells with new fill colors, and moves the active selection to U10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# New fill colors used below (BGR-packed decimal values for Range.Interior.Color)
#   10498160 = 0xA0, 0x30, 0x70 -> RGB(0x70,0x30,0xA0) purple   FF7030A0
#   15849926 = RGB(0xC6,0xD9,0xF1) "Text 2, Lighter 80%"        theme 3 / tint 0.8
#    6684927 = RGB(0xFF,0x00,0x66) pink                          FFFF0066
#    5540500 = RGB(0x94,0x8A,0x54) "Background 2, Darker 50%"   theme 2 / tint -0.5
#   16777215 = RGB(0xFF,0xFF,0xFF) white                         theme 0
#   12106214 = RGB(0xE6,0xB9,0xB8) "Accent 2, Lighter 60%"      theme 5 / tint 0.6
#    5287936 = RGB(0x00,0xB0,0x50) green                         FF00B050
#      49407 = RGB(0xFF,0xC0,0x00) gold                          FFFFC000
# ---------------------------------------------------------------------

# Recolor existing room-label marker cells (value/text unchanged, just a new fill)
$ws.Range("I1").Interior.Color  = 10498160
$ws.Range("O1").Interior.Color  = 15849926
$ws.Range("D5").Interior.Color  = 6684927
$ws.Range("I5").Interior.Color  = 5540500
$ws.Range("M6").Interior.Color  = 15849926
$ws.Range("C7").Interior.Color  = 16777215
$ws.Range("R7").Interior.Color  = 12106214
$ws.Range("W8").Interior.Color  = 10498160

# Carve the new "J" room out of the "M" room (rows 9-16, cols T:W)
$c = $ws.Range("T9");  $c.Value = "w"; $c.Interior.Color = 65535
$c = $ws.Range("U9");  $c.Value = "w"; $c.Interior.Color = 65535
$c = $ws.Range("V9");  $c.Value = "w"; $c.Interior.Color = 65535

$c = $ws.Range("T10"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U10").Value = "J"
$ws.Range("V10").Value = "J"
$ws.Range("W10").Value = "J"

$ws.Range("H11").Interior.Color = 5287936

$c = $ws.Range("T11"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U11").Value = "J"
$ws.Range("V11").Value = "J"
$ws.Range("W11").Value = "J"

$ws.Range("O12").Interior.Color = 5540500

$c = $ws.Range("T12"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U12").Value = "JL"
$ws.Range("V12").Value = "J"
$ws.Range("W12").Value = "J"

$c = $ws.Range("T13"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U13").Value = "JL"
$ws.Range("V13").Value = "J"
$ws.Range("W13").Value = "J"

$ws.Range("E14").Interior.Color = 5540500

$c = $ws.Range("T14"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U14").Value = "J"
$ws.Range("V14").Value = "J"
$ws.Range("W14").Value = "J"

$ws.Range("A15").Interior.Color = 10498160
$ws.Range("G15").Interior.Color = 15849926

$c = $ws.Range("T15"); $c.Value = "w"; $c.Interior.Color = 65535
$ws.Range("U15").Value = "J"
$ws.Range("V15").Value = "J"
$ws.Range("W15").Value = "J"

$c = $ws.Range("T16"); $c.Value = "w"; $c.Interior.Color = 65535
$c = $ws.Range("U16"); $c.Value = "w"; $c.Interior.Color = 65535
$c = $ws.Range("V16"); $c.Value = "w"; $c.Interior.Color = 65535

# More room-label marker recolors
$ws.Range("D17").Interior.Color = 16777215
$ws.Range("P17").Interior.Color = 49407
$ws.Range("S17").Interior.Color = 16777215

$ws.Range("H18").Interior.Color = 5287936
$ws.Range("L18").Interior.Color = 5540500

$ws.Range("S19").Interior.Color = 6684927

$ws.Range("H22").Interior.Color = 15849926
$ws.Range("O22").Interior.Color = 10498160

# Move the active selection, matching the saved view state in the workbook
$ws.Range("U10").Select()
